# Apply cryptos list update (prices / volume / coin rename) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.160.63'
$ws.Range('E2').Value = '  +3.65%  '
$ws.Range('D3').Value = '2.244.27'
$ws.Range('E3').Value = '  +1.31%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'294.83"
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').Value = "'86.80"
$ws.Range('E6').Value = '  +9.48%  '
$ws.Range('D7').Value = "'0.517"
$ws.Range('E7').Value = '  +2.40%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = "'0.474"
$ws.Range('E9').Value = '  +4.08%  '
$ws.Range('D10').Value = "'31.21"
$ws.Range('E10').Value = '  +12.30%  '
$ws.Range('D11').Value = "'0.0800"
$ws.Range('E11').Value = '  +3.84%  '
$ws.Range('D12').Value = "'47.25"
$ws.Range('E12').Value = '  +2.43%  '
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').Value = "'6.48"
$ws.Range('E14').Value = '  +6.59%  '
$ws.Range('D15').Value = '2.592.96'
$ws.Range('E15').Value = '  +1.07%  '
$ws.Range('D16').Value = "'14.25"
$ws.Range('D17').Value = '2.231.87'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').Value = "'0.738"
$ws.Range('E18').Value = '  +4.01%  '
$ws.Range('D19').Value = '40.095.03'
$ws.Range('E19').Value = '  +3.58%  '
$ws.Range('D20').Value = '0.0₃0895'
$ws.Range('E20').Value = '  +4.60%  '
$ws.Range('D21').Value = "'5.84"
$ws.Range('E21').Value = '  +2.24%  '
$ws.Range('D22').Value = "'10.67"
$ws.Range('E22').Value = '  +9.11%  '
$ws.Range('D23').Value = "'65.71"
$ws.Range('E23').Value = '  +1.65%  '
$ws.Range('D24').Value = "'236.64"
$ws.Range('E24').Value = '  +5.24%  '
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('D26').Value = "'2.47"
$ws.Range('E26').Value = '  +3.92%  '
$ws.Range('E27').Value = '  +8.24%  '
$ws.Range('D28').Value = "'23.00"
$ws.Range('E28').Value = '  +4.27%  '
$ws.Range('D29').Value = "'2.23"
$ws.Range('D30').Value = "'9.28"
$ws.Range('E30').Value = '  +5.00%  '
$ws.Range('D31').Value = "'33.38"
$ws.Range('E31').Value = '  +7.94%  '
$ws.Range('D32').Value = "'153.87"
$ws.Range('E32').Value = '  +3.82%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').Value = "'4.91"
$ws.Range('E34').Value = '  +3.25%  '
$ws.Range('D35').Value = "'0.0719"
$ws.Range('E35').Value = '  +5.67%  '
$ws.Range('E36').Value = '  +3.36%  '
$ws.Range('D37').Value = "'16.59"
$ws.Range('E37').Value = '  +16.06%  '
$ws.Range('E38').Value = '  +6.41%  '
$ws.Range('E39').Value = '  +3.11%  '
$ws.Range('D40').Value = "'2.72"
$ws.Range('D41').Value = "'1.70"
$ws.Range('E41').Value = '  +7.38%  '
$ws.Range('E42').Value = '  +6.78%  '
$ws.Range('D43').Value = '2.024.25'
$ws.Range('E43').Value = '  +6.63%  '
$ws.Range('E44').Value = '  +10.52%  '
$ws.Range('D45').Value = "'0.0272"
$ws.Range('E45').Value = '  +7.81%  '
$ws.Range('E46').Value = '  +11.46%  '
$ws.Range('D47').Value = "'16.34"
$ws.Range('E47').Value = '  +2.19%  '
$ws.Range('D48').Value = "'2.58"
$ws.Range('E48').Value = '  +3.47%  '
$ws.Range('D49').Value = '2.473.20'
$ws.Range('D50').Value = "'71.55"
$ws.Range('E50').Value = '  +4.96%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = "'1.12"
$ws.Range('E51').Value = '  +7.49%  '
